$d = $word.ActiveDocument

# 1. Merge the title runs "Problem: A Cat, " + "A Parrot, and a Bag of Seed"
#    into a single run, removing the _GoBack bookmark that sat between them.
$d.Content.Find.Execute("Problem: A Cat, A Parrot, and a Bag of Seed", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Problem: A Cat, A Parrot, and a Bag of Seed", 2)

# 2. Append additional solution text after the existing sentence.
$d.Content.Find.Execute("Go back and get the parrot.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Go back and get the parrot. Another solution would be the same exact scenario but taking the bag of seeds instead of the cat, and then get bring back the parrot, leave it and bring back the cat, then go back for the parrot.", 2)
